$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Brands")
$ws2 = $wb.Worksheets.Item("Special Circumstances")

$data = @(
    @('id', 'page-title', 'page-description', 'banner', 'bannerlink', 'banneralt', 'promo'),
    @('Brands-Apparel-Merchandise', 'Brands Apparel & Merchandise at CampusColors.com', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'blank', $null, 'Shop By Brand Banner', 'Some promotional text goes here'),
    @('team_long', 'page-description', 'team_league', 'team_short', 'team_code-img', 'sub_banner', 'max-price'),
    @('''47 Brand', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'BS    ', 'blank', 150),
    @('''47 Brand', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'P-TWN', 'blank', 150),
    @('''47 Brand', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'TWIN  ', 'blank', 150),
    @('Adidas', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'ADDS  ', 'blank', 150),
    @('Adidas', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'YADDS ', 'blank', 150),
    @('All-Sportz Brush', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'ASB', 'blank', 150),
    @('Aminco', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'AMIN  ', 'blank', 150),
    @('Aminco', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'P-AMIN', 'blank', 150),
    @('Antigua', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'ANT   ', 'blank', 150),
    @('Bleacher Creatures', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'BC    ', 'blank', 150),
    @('Blue 84', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'B84   ', 'blank', 150),
    @('Boelter', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'BOE   ', 'blank', 150),
    @('BSI Products', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'BSI   ', 'blank', 150),
    @('BSI Products', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'P-BSI', 'blank', 150),
    @('C & I Collectibles', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'CIC   ', 'blank', 150),
    @('Campus Colors', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'NA    ', 'blank', 150),
    @('Campus Colors', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'UT    ', 'blank', 150),
    @('Champion', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'CHMP  ', 'blank', 150),
    @('Charm 14', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'CH14  ', 'blank', 150),
    @('Colosseum', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'COL   ', 'blank', 150),
    @('Columbia', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'CSU   ', 'blank', 150),
    @('Comfy Feet', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'CMFY  ', 'blank', 150),
    @('Concept One', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'CONE  ', 'blank', 150),
    @('Concept One', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'P-CONE', 'blank', 150),
    @('Concept Sport', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'CC    ', 'blank', 150),
    @('Coopersburg Sports', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'COOP  ', 'blank', 150),
    @('Dallas Cowboys Merchandise', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'P-DCM', 'blank', 150),
    @('Distant Replays', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'DR    ', 'blank', 150),
    @('Duck House Sports', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'P-DHS', 'blank', 150),
    @('Elite Image Eyewear', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'EIE', 'blank', 150),
    @('Fabrique', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'SYK   ', 'blank', 150),
    @('Fan Mats', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'FANMAT', 'blank', 150),
    @('Fan Tape', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'FT    ', 'blank', 150),
    @('Fanatic Group', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'FAN   ', 'blank', 150),
    @('Fathead', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'FHD   ', 'blank', 150),
    @('Fit 2 Win', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'F2W   ', 'blank', 150),
    @('Foam Fanatics', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'FOAM  ', 'blank', 150),
    @('For Bare Feet', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'FBF   ', 'blank', 150),
    @('For Bare Feet', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'P-FBF', 'blank', 150),
    @('Franklin ', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'FS    ', 'blank', 150),
    @('Freemont Die', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'FRE   ', 'blank', 150),
    @('Front Row Sports', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'FRS   ', 'blank', 150),
    @('Game Master', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'GM    ', 'blank', 150),
    @('Gameday Spirit', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'GS', 'blank', 150),
    @('Gear', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'GEAR  ', 'blank', 150),
    @('Genuine Stuff', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'OS    ', 'blank', 150),
    @('GIII', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'GIII  ', 'blank', 150),
    @('Glitter Girl', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'GLIT', 'blank', 150),
    @('Glory Haus', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'GH    ', 'blank', 150),
    @('Great American Products', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'GAP   ', 'blank', 150),
    @('Highland Mint', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'HM    ', 'blank', 150),
    @('Highland Mint', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'THM   ', 'blank', 150),
    @('Holloway', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'HOL   ', 'blank', 150),
    @('Hunter', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'HUN   ', 'blank', 150),
    @('Hunter', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'P-HUN', 'blank', 150),
    @('J America', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'JAM   ', 'blank', 150),
    @('Jarden Sports', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'JSL   ', 'blank', 150),
    @('Kasky Kids', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'KAS   ', 'blank', 150),
    @('League Collegiate Outfitters', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'LCO', 'blank', 150),
    @('Letter Art', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'PCH   ', 'blank', 150),
    @('Little Earth', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'LE    ', 'blank', 150),
    @('Logo', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'LOGO  ', 'blank', 150),
    @('Majestic', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'MAJ   ', 'blank', 150),
    @('Majestic', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'P-MAJ', 'blank', 150),
    @('Majestic', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'YMAJ  ', 'blank', 150),
    @('Mascot Factory', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'MAS', 'blank', 150),
    @('Me & My Big Ideas', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'MMBI', 'blank', 150),
    @('Mojo', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'RI    ', 'blank', 150),
    @('MV Sport', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'MV', 'blank', 150),
    @('New Era', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'NE    ', 'blank', 150),
    @('New Era', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'P-NE', 'blank', 150),
    @('NFL', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'YNFL  ', 'blank', 150),
    @('Nike', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'NIKE  ', 'blank', 150),
    @('Nike', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'NIKE-Y', 'blank', 150),
    @('Original Retro Brand', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'RB    ', 'blank', 150),
    @('Oyo Sports', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'OYO   ', 'blank', 150),
    @('Photo File', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'PF    ', 'blank', 150),
    @('Pinemeadow Green', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'PINE', 'blank', 150),
    @('PPW Toys', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'PPW   ', 'blank', 150),
    @('Pro Specialties Group', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'PSG   ', 'blank', 150),
    @('Pubs Of', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'PUB   ', 'blank', 150),
    @('R & R Imports', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'RR    ', 'blank', 150),
    @('Rawlings', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'K2    ', 'blank', 150),
    @('Reebok', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'REEB  ', 'blank', 150),
    @('Reebok', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'YREEB ', 'blank', 150),
    @('Rico', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'P-RICO', 'blank', 150),
    @('Rico', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'RICO  ', 'blank', 150),
    @('Ridell', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'P-RID', 'blank', 150),
    @('Ridell', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'RID   ', 'blank', 150),
    @('Russel', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'RUSS  ', 'blank', 150),
    @('Sikiyou ', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'SISK  ', 'blank', 150),
    @('Silver Star', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'SLVS  ', 'blank', 150),
    @('Skootz', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'SKTZ  ', 'blank', 150),
    @('Soffe', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'SOFE  ', 'blank', 150),
    @('Spalding', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'SPLDG ', 'blank', 150),
    @('Sports Coverage', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'SPCOV ', 'blank', 150),
    @('Stockdale', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'STCK  ', 'blank', 150),
    @('Sutter''s Mill', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'SUT', 'blank', 150),
    @('Tailgate Clothing Co.', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'TAIL  ', 'blank', 150),
    @('Team Edition', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'TMED  ', 'blank', 150),
    @('Team Effort', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'TE    ', 'blank', 150),
    @('Team Golf', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'TG    ', 'blank', 150),
    @('Team ProMark', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'PROMRK', 'blank', 150),
    @('Team Sports America', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'P-TSA', 'blank', 150),
    @('Tervis', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'TER   ', 'blank', 150),
    @('The Emblem Source', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'TES   ', 'blank', 150),
    @('The Game', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'GAME  ', 'blank', 150),
    @('The Memory Company', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'TMC   ', 'blank', 150),
    @('The Northwest', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'NW    ', 'blank', 150),
    @('The Northwest', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'P-NW', 'blank', 150),
    @('Top of the World', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'TOW   ', 'blank', 150),
    @('Top Sox', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'TPSX  ', 'blank', 150),
    @('Topperscot', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'TOP   ', 'blank', 150),
    @('Topperscot', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'P-TOP', 'blank', 150),
    @('Under Armour', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'UA    ', 'blank', 150),
    @('Wes & Willy ', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'WW    ', 'blank', 150),
    @('Wincraft', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'MCAR  ', 'blank', 150),
    @('Wincraft', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'P-MCAR', 'blank', 150),
    @('Wincraft', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'P-WIN', 'blank', 150),
    @('Wincraft', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'WIN   ', 'blank', 150),
    @('Winning Streak Sports', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'P-WSS', 'blank', 150),
    @('Winning Streak Sports', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'WSS   ', 'blank', 150),
    @('Zephyr', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'ZEPH  ', 'blank', 150),
    @('Zipway', 'Shop Campus Colors for thousands of NCAA, NFL, NBA, MLB, & NHL products, novelties and more! We offer gear from top brands such as Nike & Adidas. Ship Same-Day to All 50 States!', 'vendor_code', 'blank', 'ZIP   ', 'blank', 150)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 1
    for ($j = 0; $j -lt 7; $j++) {
        $val = $row[$j]
        if ($null -ne $val) {
            $ws1.Cells.Item($r, $j + 1).Value = $val
        }
    }
}

$ws2.Range('A1').Value = 'Several vendors have multiple codes. Until these can be accounted for in the utility they have been consolodated by hand in the json.  These vendors include (but may not be limited to): ''47 Brand, Adidas, Aminco, BSI Products, Campus Colors, Concept One, For Bare Feet, Highland Mint, Hunter, Majestic, New Era, Nike, Reebok, Rico, Ridell, The Northwest, Topperscot, Wincraft, Winning Streak Sports.'
$ws2.Range('A2').Value = 'The id for 47 Brand includes a singe quote: ''47 Brand. This breaks elastic for this category, so remove it if the data is re-populated in the future. '
